$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Stevie Damrel's time from "16h 30m" to "13h 30m"
$ws.Range("B8").Value = "13h 30m"

# Update Maggie Burton's time from "16h 30m" to "19h 30m"
$ws.Range("B7").Value = "19h 30m"

# Move active cell selection to B7 (as reflected in the saved file)
$ws.Range("B7").Select()
